$d = $word.ActiveDocument

# Helper: replace the run content of an existing paragraph (keeps the
# paragraph's own <w:pPr> / rsid attributes untouched) with the runs
# described by $runsXml (one or more <w:r>...</w:r> elements).
function Replace-ParaRuns($para, $runsXml) {
    $r = $para.Range
    $textRange = $d.Range($r.Start, $r.End - 1)
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $textRange.InsertXML($xml)
}

# --- 1) "Uso de etiqueta obsoleta center, modificación a través de css."
#        -> split into 3 runs, and "css" -> "CSS"
$p6 = $d.Paragraphs.Item(6)
$runs6 = '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Uso de etiqueta obsoleta cent</w:t></w:r>' + `
         '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>er, modificación a través de CSS</w:t></w:r>' + `
         '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r>'
Replace-ParaRuns $p6 $runs6

# --- 2) "Etiqueta main mal cerrada en sección CuartoESO."
#        -> split into 2 runs, appending "html."
$p7 = $d.Paragraphs.Item(7)
$runs7 = '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Etiqueta main mal cerrada en sección CuartoESO.</w:t></w:r>' + `
         '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>html.</w:t></w:r>'
Replace-ParaRuns $p7 $runs7

# --- 3) Append two brand-new numbered-list paragraphs after paragraph 7
$p7 = $d.Paragraphs.Item(7)
$insPos = $p7.Range.End - 1
$collapsed = $d.Range($insPos, $insPos)

$pPrCommon = '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>'

$para8 = '<w:p>' + $pPrCommon + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Faltas de etiquetas formulario en la página Solicitud.html</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> y también en esta página había que encerrar los input en etiquetas label</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r>' + `
    '</w:p>'

$para9 = '<w:p>' + $pPrCommon + `
    '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Cambios de los alt de los logotipos que aparecen en la cabecera y en los pies de página.</w:t></w:r>' + `
    '</w:p>'

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $para8 + $para9 + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$collapsed.InsertXML($xml)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "$i => [$($p.Range.Text)]"
}
